# Scheduled market-data refresh: update currentAveragePrice* / Leve profit
# columns (H:N) on each job sheet of the Ixion_Profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 5406.3335
$ws.Range("I16").Value = 5004.5
$ws.Range("J16").Value = 6210
$ws.Range("K16").Value = 5004.5
$ws.Range("L16").Value = 6210
$ws.Range("M16").Value = -4774.5
$ws.Range("N16").Value = -6670
$ws.Range("H125").Value = 745.3333
$ws.Range("I125").Value = 400
$ws.Range("J125").Value = 918
$ws.Range("K125").Value = 3600
$ws.Range("L125").Value = 8262
$ws.Range("M125").Value = -1140
$ws.Range("N125").Value = -13182
$ws.Range("H129").Value = 1111.8387
$ws.Range("I129").Value = 619.25
$ws.Range("J129").Value = 1145.8103
$ws.Range("K129").Value = 1857.75
$ws.Range("L129").Value = 3437.4309
$ws.Range("M129").Value = 3142.25
$ws.Range("N129").Value = -13437.4309
$ws.Range("H132").Value = 2284.311
$ws.Range("I132").Value = 1569.2572
$ws.Range("J132").Value = 4787
$ws.Range("K132").Value = 4707.7716
$ws.Range("L132").Value = 14361
$ws.Range("M132").Value = -2177.7716
$ws.Range("N132").Value = -19421
$ws.Range("H137").Value = 1536.0358
$ws.Range("I137").Value = 1125.0444
$ws.Range("J137").Value = 3217.3635
$ws.Range("K137").Value = 3375.1332
$ws.Range("L137").Value = 9652.0905
$ws.Range("M137").Value = -825.1332000000002
$ws.Range("N137").Value = -14752.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2874.28
$ws.Range("I32").Value = 2874.28
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2874.28
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2587.28
$ws.Range("N32").Value = ""
$ws.Range("H33").Value = 5965.4
$ws.Range("I33").Value = 913.5
$ws.Range("J33").Value = 9333.333000000001
$ws.Range("K33").Value = 913.5
$ws.Range("L33").Value = 9333.333000000001
$ws.Range("M33").Value = -584.5
$ws.Range("N33").Value = -9991.333000000001
$ws.Range("H61").Value = 339415.5
$ws.Range("I61").Value = 6802.2173
$ws.Range("J61").Value = 1432287.8
$ws.Range("K61").Value = 6802.2173
$ws.Range("L61").Value = 1432287.8
$ws.Range("M61").Value = -6590.2173
$ws.Range("N61").Value = -1432711.8
$ws.Range("H74").Value = 1735.725
$ws.Range("I74").Value = 1398.0952
$ws.Range("J74").Value = 2108.8948
$ws.Range("K74").Value = 1398.0952
$ws.Range("L74").Value = 2108.8948
$ws.Range("M74").Value = -524.0952
$ws.Range("N74").Value = -3856.8948
$ws.Range("H77").Value = 1735.725
$ws.Range("I77").Value = 1398.0952
$ws.Range("J77").Value = 2108.8948
$ws.Range("K77").Value = 6990.476
$ws.Range("L77").Value = 10544.474
$ws.Range("M77").Value = -2622.476
$ws.Range("N77").Value = -19280.474
$ws.Range("H97").Value = 1007.913
$ws.Range("I97").Value = 1109.15
$ws.Range("J97").Value = 333
$ws.Range("K97").Value = 1109.15
$ws.Range("L97").Value = 333
$ws.Range("M97").Value = -613.1500000000001
$ws.Range("N97").Value = -1325
$ws.Range("H122").Value = 1711639.5
$ws.Range("I122").Value = 1833778
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 5501334
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -5498884
$ws.Range("N122").Value = -10000
$ws.Range("H132").Value = 1494835.9
$ws.Range("I132").Value = 1723.3405
$ws.Range("J132").Value = 5003650.5
$ws.Range("K132").Value = 5170.0215
$ws.Range("L132").Value = 15010951.5
$ws.Range("M132").Value = -2640.0215
$ws.Range("N132").Value = -15016011.5
$ws.Range("H136").Value = 339415.5
$ws.Range("I136").Value = 6802.2173
$ws.Range("J136").Value = 1432287.8
$ws.Range("K136").Value = 20406.6519
$ws.Range("L136").Value = 4296863.4
$ws.Range("M136").Value = -17856.6519
$ws.Range("N136").Value = -4301963.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1766.7
$ws.Range("I94").Value = 1224.9231
$ws.Range("J94").Value = 2772.8572
$ws.Range("K94").Value = 1224.9231
$ws.Range("L94").Value = 2772.8572
$ws.Range("M94").Value = -773.9231
$ws.Range("N94").Value = -3674.8572
$ws.Range("H107").Value = 1375.421
$ws.Range("I107").Value = 1383.75
$ws.Range("J107").Value = 1331
$ws.Range("K107").Value = 1383.75
$ws.Range("L107").Value = 1331
$ws.Range("M107").Value = 536.25
$ws.Range("N107").Value = -5171
$ws.Range("H134").Value = 33884.914
$ws.Range("I134").Value = 7528.8945
$ws.Range("J134").Value = 65182.688
$ws.Range("K134").Value = 22586.6835
$ws.Range("L134").Value = 195548.064
$ws.Range("M134").Value = -20051.6835
$ws.Range("N134").Value = -200618.064

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2294.077
$ws.Range("I16").Value = 2077.75
$ws.Range("J16").Value = 2390.2222
$ws.Range("K16").Value = 2077.75
$ws.Range("L16").Value = 2390.2222
$ws.Range("M16").Value = -1790.75
$ws.Range("N16").Value = -2964.2222
$ws.Range("H31").Value = 4593.9414
$ws.Range("I31").Value = 1962.7046
$ws.Range("J31").Value = 9417.875
$ws.Range("K31").Value = 1962.7046
$ws.Range("L31").Value = 9417.875
$ws.Range("M31").Value = -1667.7046
$ws.Range("N31").Value = -10007.875
$ws.Range("H34").Value = 4593.9414
$ws.Range("I34").Value = 1962.7046
$ws.Range("J34").Value = 9417.875
$ws.Range("K34").Value = 1962.7046
$ws.Range("L34").Value = 9417.875
$ws.Range("M34").Value = -1760.7046
$ws.Range("N34").Value = -9821.875
$ws.Range("H105").Value = 3119.4119
$ws.Range("I105").Value = 3126.875
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 3126.875
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -1379.875
$ws.Range("N105").Value = -6494
$ws.Range("H113").Value = 2294.077
$ws.Range("I113").Value = 2077.75
$ws.Range("J113").Value = 2390.2222
$ws.Range("K113").Value = 2077.75
$ws.Range("L113").Value = 2390.2222
$ws.Range("M113").Value = 92.25
$ws.Range("N113").Value = -6730.2222
$ws.Range("H122").Value = 3059.4
$ws.Range("I122").Value = 3177.111
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 9531.332999999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -7081.332999999999
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 2506.5557
$ws.Range("I132").Value = 1299.1111
$ws.Range("J132").Value = 3714
$ws.Range("K132").Value = 3897.3333
$ws.Range("L132").Value = 11142
$ws.Range("M132").Value = -1367.3333
$ws.Range("N132").Value = -16202
$ws.Range("H134").Value = 178783.72
$ws.Range("I134").Value = 3623.25
$ws.Range("J134").Value = 479058.8
$ws.Range("K134").Value = 10869.75
$ws.Range("L134").Value = 1437176.4
$ws.Range("M134").Value = -8334.75
$ws.Range("N134").Value = -1442246.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 285.81818
$ws.Range("I6").Value = 60.333332
$ws.Range("J6").Value = 370.375
$ws.Range("K6").Value = 180.999996
$ws.Range("L6").Value = 1111.125
$ws.Range("M6").Value = -67.99999600000001
$ws.Range("N6").Value = -1337.125
$ws.Range("H10").Value = 173.5
$ws.Range("I10").Value = 106
$ws.Range("J10").Value = 274.75
$ws.Range("K10").Value = 318
$ws.Range("L10").Value = 824.25
$ws.Range("M10").Value = -179
$ws.Range("N10").Value = -1102.25
$ws.Range("H13").Value = 2607.6924
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 2990.9092
$ws.Range("K13").Value = 1500
$ws.Range("L13").Value = 8972.7276
$ws.Range("M13").Value = -1332
$ws.Range("N13").Value = -9308.7276
$ws.Range("H21").Value = 1272.6364
$ws.Range("I21").Value = 500
$ws.Range("J21").Value = 1349.9
$ws.Range("K21").Value = 1500
$ws.Range("L21").Value = 4049.7
$ws.Range("M21").Value = -1327
$ws.Range("N21").Value = -4395.700000000001
$ws.Range("H26").Value = 460.92307
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 491
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 1473
$ws.Range("M26").Value = -12
$ws.Range("N26").Value = -2049
$ws.Range("H34").Value = 1520
$ws.Range("I34").Value = 850
$ws.Range("J34").Value = 1698.6666
$ws.Range("K34").Value = 2550
$ws.Range("L34").Value = 5095.9998
$ws.Range("M34").Value = -2466
$ws.Range("N34").Value = -5263.9998
$ws.Range("H55").Value = 3442.8572
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3442.8572
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10328.5716
$ws.Range("N55").Value = -10682.5716
$ws.Range("H123").Value = 8608.333000000001
$ws.Range("I123").Value = 7850
$ws.Range("J123").Value = 8987.5
$ws.Range("K123").Value = 23550
$ws.Range("L123").Value = 26962.5
$ws.Range("M123").Value = -21100
$ws.Range("N123").Value = -31862.5
$ws.Range("H131").Value = 2440022.8
$ws.Range("I131").Value = 9091575
$ws.Range("J131").Value = 1120.2333
$ws.Range("K131").Value = 27274725
$ws.Range("L131").Value = 3360.699900000001
$ws.Range("M131").Value = -27269685
$ws.Range("N131").Value = -13440.6999
$ws.Range("H139").Value = 4017.4255
$ws.Range("I139").Value = 4896.76
$ws.Range("J139").Value = 3018.182
$ws.Range("K139").Value = 14690.28
$ws.Range("L139").Value = 9054.545999999998
$ws.Range("M139").Value = -9550.280000000001
$ws.Range("N139").Value = -19334.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 11900
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 11900
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 11900
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -12476
$ws.Range("H21").Value = 29333.334
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 29333.334
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 29333.334
$ws.Range("N21").Value = -29679.334
$ws.Range("H30").Value = 29333.334
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 29333.334
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 29333.334
$ws.Range("N30").Value = -29543.334
$ws.Range("H53").Value = 19800
$ws.Range("I53").Value = 19500
$ws.Range("J53").Value = 19900
$ws.Range("K53").Value = 19500
$ws.Range("L53").Value = 19900
$ws.Range("M53").Value = -18869
$ws.Range("N53").Value = -21162
$ws.Range("H80").Value = 6359.36
$ws.Range("I80").Value = 8555.3125
$ws.Range("J80").Value = 2455.4443
$ws.Range("K80").Value = 8555.3125
$ws.Range("L80").Value = 2455.4443
$ws.Range("M80").Value = -7557.3125
$ws.Range("N80").Value = -4451.4443
$ws.Range("H83").Value = 6359.36
$ws.Range("I83").Value = 8555.3125
$ws.Range("J83").Value = 2455.4443
$ws.Range("K83").Value = 42776.5625
$ws.Range("L83").Value = 12277.2215
$ws.Range("M83").Value = -37784.5625
$ws.Range("N83").Value = -22261.2215
$ws.Range("H97").Value = 1824.5294
$ws.Range("I97").Value = 1709
$ws.Range("J97").Value = 2101.8
$ws.Range("K97").Value = 1709
$ws.Range("L97").Value = 2101.8
$ws.Range("M97").Value = -1213
$ws.Range("N97").Value = -3093.8
$ws.Range("H122").Value = 106484000
$ws.Range("I122").Value = 106484000
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 319452000
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -319449550
$ws.Range("H123").Value = 14225.556
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 14225.556
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 14225.556
$ws.Range("N123").Value = -19125.556
$ws.Range("H132").Value = 5820.1313
$ws.Range("I132").Value = 6689.9614
$ws.Range("J132").Value = 3935.5
$ws.Range("K132").Value = 20069.8842
$ws.Range("L132").Value = 11806.5
$ws.Range("M132").Value = -17539.8842
$ws.Range("N132").Value = -16866.5
$ws.Range("H136").Value = 10929.333
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 10929.333
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 32787.999
$ws.Range("N136").Value = -37887.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 150000
$ws.Range("I25").Value = 150000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 150000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -149770
$ws.Range("N25").Value = ""
$ws.Range("H93").Value = 1669.25
$ws.Range("I93").Value = 1647.2
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1647.2
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -399.2
$ws.Range("N93").Value = -4496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = ""
